$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the bollinger band trade calculation:
# - fill in the SellPrice (E2) that was missing
# - compute the Price Change % (F2) between BuyPrice and SellPrice
# - correct the Profitable (B2) / Holding (G2) flags now that the trade has closed
# - add row 3 with the updated Principle (C3) after the trade closed
$ws.Range("B2").Value = $true
$ws.Range("E2").Value = 108.91
$ws.Range("F2").Value = 0.45194613539936812
$ws.Range("G2").Value = $false

$ws.Range("C3").Value = 10045.19

# Column C widens slightly to fit the new Principle value
$ws.Columns.Item(3).ColumnWidth = 8.14
